# "Loading all tables complete" -- finish describing the "answers" sheet in
# the metadata table: fix up the description text for the answers table,
# wrap it so it's readable in the cell, and resize the rows to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# B2 held a description that used to start with "This table contains ...".
# Re-word it to start with the sheet name it is actually describing.
$ws.Range("B2").Value = "answers table contains the results of of an exam. The first row contains the index for the questions, and the following column contains the answer to each question for each student. 1 means they provided the right answer, 0 means they provided the wrong answer"

# Wrap the text in column B so the long description is readable, and grow
# the rows to fit the wrapped content.
$ws.Columns.Item(2).WrapText = $true
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 85

# Leave the selection on A2.
$ws.Range("A2").Select()
